$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text formatting (not auto-converted to numbers)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.751.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.772.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.09"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.770.28"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.85%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.14%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.02%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.443"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.19%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.68"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +6.13%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.61%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.405.72"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.773.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.771.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "457.41"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.77%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.46"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.49%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.18%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.69"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.71%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.82"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.08%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.77%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.918.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.43%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.17"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.78%  "

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -6.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.83"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.94"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.143"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.37%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.976"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.37%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.16"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Arweave"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.57"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "152.22"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.14%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.07%  "

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.27"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "386.02"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.44%  "
